# Mark slots with incomplete info (missing time interval or people needed)
# so they can be skipped downstream.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# WED row 5 ("people needed" slot) was "1*" -> now a placeholder "dooo"
$ws.Range("G5").Value = "dooo"

# TUE row 6 ("time interval" slot) was "11:30-14:30" -> now a placeholder "hello"
$ws.Range("E6").Value = "hello"
